# Insert two new price-report rows (Florida King, Primera/Segunda) at the
# top of the "Durazno" block for Terminal Hortofrutícola Agro Chillán,
# pushing the existing rows 232-291 down to 234-293.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("232:233").Insert()

# New row 232: Florida King / Primera
$ws.Range("A232").Value2 = 7
$ws.Range("B232").Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C232").Value2 = "Ñuble"
$ws.Range("D232").Value2 = 44900
$ws.Range("E232").Value2 = 16
$ws.Range("F232").Value2 = "Fruta"
$ws.Range("G232").Value2 = 100103
$ws.Range("H232").Value2 = "Frutos de hueso (carozo)"
$ws.Range("I232").Value2 = 100103004
$ws.Range("J232").Value2 = "Durazno"
$ws.Range("K232").Value2 = "Florida King"
$ws.Range("L232").Value2 = "Primera"
$ws.Range("M232").Value2 = 160
$ws.Range("N232").Value2 = 11000
$ws.Range("O232").Value2 = 12000
$ws.Range("P232").Value2 = 11500
$ws.Range("Q232").Value2 = "$/caja 15 kilos granel"
$ws.Range("R232").Value2 = "Región de O'Higgins"
$ws.Range("S232").Value2 = 767
$ws.Range("T232").Value2 = 15

# New row 233: Florida King / Segunda
$ws.Range("A233").Value2 = 7
$ws.Range("B233").Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C233").Value2 = "Ñuble"
$ws.Range("D233").Value2 = 44900
$ws.Range("E233").Value2 = 16
$ws.Range("F233").Value2 = "Fruta"
$ws.Range("G233").Value2 = 100103
$ws.Range("H233").Value2 = "Frutos de hueso (carozo)"
$ws.Range("I233").Value2 = 100103004
$ws.Range("J233").Value2 = "Durazno"
$ws.Range("K233").Value2 = "Florida King"
$ws.Range("L233").Value2 = "Segunda"
$ws.Range("M233").Value2 = 80
$ws.Range("N233").Value2 = 10000
$ws.Range("O233").Value2 = 10000
$ws.Range("P233").Value2 = 10000
$ws.Range("Q233").Value2 = "$/caja 15 kilos granel"
$ws.Range("R233").Value2 = "Región de O'Higgins"
$ws.Range("S233").Value2 = 667
$ws.Range("T233").Value2 = 15
